$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(432, "2021-06-11", 76, 74),
  @(433, "2021-06-12", 55, 52),
  @(434, "2021-06-13", 63, 63),
  @(435, "2021-06-14", 68, 63),
  @(436, "2021-06-15", 81, 73),
  @(437, "2021-06-16", 76, 69),
  @(438, "2021-06-17", 74, 71),
  @(439, "2021-06-18", 73, 70),
  @(440, "2021-06-19", 64, 59),
  @(441, "2021-06-20", 72, 69),
  @(442, "2021-06-21", 72, 67),
  @(443, "2021-06-22", 68, 66),
  @(444, "2021-06-23", 72, 67),
  @(445, "2021-06-24", 81, 73),
  @(446, "2021-06-25", 83, 72),
  @(447, "2021-06-26", 61, 56),
  @(448, "2021-06-27", 73, 69),
  @(449, "2021-06-28", 67, 66),
  @(450, "2021-06-29", 68, 66),
  @(451, "2021-06-30", 72, 69),
  @(452, "2021-07-01", 90, 76),
  @(453, "2021-07-02", 73, 69),
)

foreach ($item in $data) {
  $r = $item[0]
  $date = $item[1]
  $b = $item[2]
  $c = $item[3]

  # Copy formatting from the row above so styles (s indices) match the existing table
  $ws.Range("A" + ($r - 1) + ":D" + ($r - 1)).Copy()
  $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

  $ws.Cells.Item($r, 1).Value = $date
  $ws.Cells.Item($r, 2).Value = $b
  $ws.Cells.Item($r, 3).Value = $c
  $ws.Cells.Item($r, 4).Formula = "=C" + $r + "/B" + $r
}

$excel.CutCopyMode = 0

$ws.Range("D437:D453").Select()
